$wb = $excel.ActiveWorkbook

# Update the Custid value (column L, row 2) on the NewCust, DeleteCust
# and EditCust sheets from "38108" to "83364" for the testng reports changes.
$sheetNames = @("NewCust", "DeleteCust", "EditCust")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("L2").Value = "83364"
}
